# Changed circumference algorithm; Up dated data and images; Added cirfumference visualisation
#
# Populate the "Circumference" column (C) for every data row (2-119) with the
# newly computed values from the updated circumference algorithm. Column C
# already has its header ("Circumference") in row 1; only the per-grain
# values are being added here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 1086.004175901413
$ws.Cells.Item(3, 3).Value = 214.3675310611725
$ws.Cells.Item(4, 3).Value = 909.3178842067719
$ws.Cells.Item(5, 3).Value = 675.5533845424652
$ws.Cells.Item(6, 3).Value = 279.4213538169861
$ws.Cells.Item(7, 3).Value = 606.6833273172379
$ws.Cells.Item(8, 3).Value = 395.1025931835175
$ws.Cells.Item(9, 3).Value = 1385.716940760612
$ws.Cells.Item(10, 3).Value = 1327.859076142311
$ws.Cells.Item(11, 3).Value = 1917.706836104393
$ws.Cells.Item(12, 3).Value = 941.903670668602
$ws.Cells.Item(13, 3).Value = 123.012192606926
$ws.Cells.Item(14, 3).Value = 678.4823168516159
$ws.Cells.Item(15, 3).Value = 53.69848430156708
$ws.Cells.Item(16, 3).Value = 67.35533845424652
$ws.Cells.Item(17, 3).Value = 75.25483322143555
$ws.Cells.Item(18, 3).Value = 226.1076455116272
$ws.Cells.Item(19, 3).Value = 296.3919162750244
$ws.Cells.Item(20, 3).Value = 177.622364282608
$ws.Cells.Item(21, 3).Value = 191.3380935192108
$ws.Cells.Item(22, 3).Value = 1566.243852138519
$ws.Cells.Item(23, 3).Value = 455.9726504087448
$ws.Cells.Item(24, 3).Value = 226.0071402788162
$ws.Cells.Item(25, 3).Value = 374.2741661071777
$ws.Cells.Item(26, 3).Value = 677.5533845424652
$ws.Cells.Item(27, 3).Value = 579.3107439279556
$ws.Cells.Item(28, 3).Value = 581.4112491607666
$ws.Cells.Item(29, 3).Value = 965.1412589550018
$ws.Cells.Item(30, 3).Value = 163.2375882863998
$ws.Cells.Item(31, 3).Value = 244.8771975040436
$ws.Cells.Item(32, 3).Value = 479.9482651948929
$ws.Cells.Item(33, 3).Value = 322.7766922712326
$ws.Cells.Item(34, 3).Value = 284.634556889534
$ws.Cells.Item(35, 3).Value = 319.1198381185532
$ws.Cells.Item(36, 3).Value = 198.4507913589478
$ws.Cells.Item(37, 3).Value = 826.7493426799774
$ws.Cells.Item(38, 3).Value = 203.8650048971176
$ws.Cells.Item(39, 3).Value = 114.5685415267944
$ws.Cells.Item(40, 3).Value = 237.0365778207779
$ws.Cells.Item(41, 3).Value = 220.7939372062683
$ws.Cells.Item(42, 3).Value = 470.9432128667831
$ws.Cells.Item(43, 3).Value = 96.56854152679443
$ws.Cells.Item(44, 3).Value = 60.87005722522736
$ws.Cells.Item(45, 3).Value = 495.9137753248215
$ws.Cells.Item(46, 3).Value = 109.6396092176437
$ws.Cells.Item(47, 3).Value = 276.1492756605148
$ws.Cells.Item(48, 3).Value = 217.0365778207779
$ws.Cells.Item(49, 3).Value = 58.18376553058624
$ws.Cells.Item(50, 3).Value = 355.64674949646
$ws.Cells.Item(51, 3).Value = 143.923879981041
$ws.Cells.Item(52, 3).Value = 475.2863587141037
$ws.Cells.Item(53, 3).Value = 600.180801153183
$ws.Cells.Item(54, 3).Value = 820.5483322143555
$ws.Cells.Item(55, 3).Value = 71.84061968326569
$ws.Cells.Item(56, 3).Value = 125.2964633703232
$ws.Cells.Item(57, 3).Value = 383.6883796453476
$ws.Cells.Item(58, 3).Value = 243.7056245803833
$ws.Cells.Item(59, 3).Value = 296.0487704277039
$ws.Cells.Item(60, 3).Value = 297.3624787330627
$ws.Cells.Item(61, 3).Value = 1068.864128470421
$ws.Cells.Item(62, 3).Value = 556.7422024011612
$ws.Cells.Item(63, 3).Value = 482.0142805576324
$ws.Cells.Item(64, 3).Value = 237.9066350460052
$ws.Cells.Item(65, 3).Value = 659.5950146913528
$ws.Cells.Item(66, 3).Value = 329.2619735002518
$ws.Cells.Item(67, 3).Value = 334.1320307254791
$ws.Cells.Item(68, 3).Value = 281.2203433513641
$ws.Cells.Item(69, 3).Value = 429.7888848781586
$ws.Cells.Item(70, 3).Value = 212.3502861261368
$ws.Cells.Item(71, 3).Value = 458.6000670194626
$ws.Cells.Item(72, 3).Value = 508.8843377828598
$ws.Cells.Item(73, 3).Value = 253.8061298131943
$ws.Cells.Item(74, 3).Value = 156.2670258283615
$ws.Cells.Item(75, 3).Value = 293.9482651948929
$ws.Cells.Item(76, 3).Value = 219.0365778207779
$ws.Cells.Item(77, 3).Value = 201.9655101299286
$ws.Cells.Item(78, 3).Value = 182.4091612100601
$ws.Cells.Item(79, 3).Value = 260.1492756605148
$ws.Cells.Item(80, 3).Value = 532.3401814699173
$ws.Cells.Item(81, 3).Value = 865.5432798862457
$ws.Cells.Item(82, 3).Value = 1108.454967260361
$ws.Cells.Item(83, 3).Value = 240.4924215078354
$ws.Cells.Item(84, 3).Value = 503.7716399431229
$ws.Cells.Item(85, 3).Value = 598.0975408554077
$ws.Cells.Item(86, 3).Value = 381.8893901109695
$ws.Cells.Item(87, 3).Value = 703.1513636112213
$ws.Cells.Item(88, 3).Value = 209.82337474823
$ws.Cells.Item(89, 3).Value = 553.9137753248215
$ws.Cells.Item(90, 3).Value = 693.4945094585419
$ws.Cells.Item(91, 3).Value = 300.2914110422134
$ws.Cells.Item(92, 3).Value = 248.5929267406464
$ws.Cells.Item(93, 3).Value = 207.8650048971176
$ws.Cells.Item(94, 3).Value = 286.3919162750244
$ws.Cells.Item(95, 3).Value = 693.2518688440323
$ws.Cells.Item(96, 3).Value = 145.5807341337204
$ws.Cells.Item(97, 3).Value = 299.0193328857422
$ws.Cells.Item(98, 3).Value = 866.0874361991882
$ws.Cells.Item(99, 3).Value = 382.2325359582901
$ws.Cells.Item(100, 3).Value = 586.2396762371063
$ws.Cells.Item(101, 3).Value = 311.9482651948929
$ws.Cells.Item(102, 3).Value = 456.7594473361969
$ws.Cells.Item(103, 3).Value = 384.4507913589478
$ws.Cells.Item(104, 3).Value = 194.4507913589478
$ws.Cells.Item(105, 3).Value = 499.0437180995941
$ws.Cells.Item(106, 3).Value = 174.7523070573807
$ws.Cells.Item(107, 3).Value = 530.3990565538406
$ws.Cells.Item(108, 3).Value = 650.9675980806351
$ws.Cells.Item(109, 3).Value = 281.5634891986847
$ws.Cells.Item(110, 3).Value = 413.8721451759338
$ws.Cells.Item(111, 3).Value = 56.87005722522736
$ws.Cells.Item(112, 3).Value = 220.1665205955505
$ws.Cells.Item(113, 3).Value = 243.7644996643066
$ws.Cells.Item(114, 3).Value = 171.82337474823
$ws.Cells.Item(115, 3).Value = 263.4213538169861
$ws.Cells.Item(116, 3).Value = 203.4802289009094
$ws.Cells.Item(117, 3).Value = 123.9827550649643
$ws.Cells.Item(118, 3).Value = 74.04163014888763
$ws.Cells.Item(119, 3).Value = 160.0243852138519
